# "Final changes for testing part of REF"
# Populate the Result sheet with the outcome of running the Tests sheet:
#  - copy WorkflowFile/Outcome rows from Tests into Result
#  - rename Result's second header from ExpectedResult to ActualResult
#  - drop the old Status/Comments columns, the AutoFilter and the
#    corresponding _FilterDatabase defined name
#  - make Result the active/selected sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Tests"
$ws2 = $wb.Worksheets.Item(2)   # "Result"

# Drop the AutoFilter on Result and the hidden _FilterDatabase name it created
if ($ws2.AutoFilterMode) {
    $ws2.AutoFilterMode = $false
}
foreach ($n in $wb.Names) {
    $n.Delete()
}

# Result no longer needs the Status / Comments columns
$ws2.Columns.Item(4).Delete()
$ws2.Columns.Item(3).Delete()

# Header: keep WorkflowFile, rename second column to ActualResult
$ws2.Range("A1").Value = "WorkflowFile"
$ws2.Range("B1").Value = "ActualResult"

# Copy the 8 test rows straight across from Tests
for ($r = 2; $r -le 9; $r++) {
    $ws2.Cells.Item($r, 1).Value = $ws1.Cells.Item($r, 1).Text
    $ws2.Cells.Item($r, 2).Value = $ws1.Cells.Item($r, 2).Text
}

# Match the bold/16pt header formatting used on the Tests sheet
$ws2.Range("A1:B1").Font.Bold = $true
$ws2.Range("A1:B1").Font.Size = 16
$ws2.Rows.Item(1).RowHeight = 21

# Re-scope the dropdown validation on column B to just the data rows
$ws2.Range("B2:B1048576").Validation.Delete()
$ws2.Range("B2:B9").Validation.Add(3, 1, 1, """Success,BusinessException,SystemException""")

# Result becomes the active sheet/tab, with the same selections seen in Excel
$ws2.Select()
$ws2.Range("G8").Select()
$ws1.Range("A1:B9").Select()
$ws2.Select()
